# Auto-update draw results: append the 2025-11-25 "Pick 4" draw as row 70.
#
# The sheet stores every column as plain text, even columns that look like
# dates/numbers (e.g. "2025-11-25", "251125"). Force Text formatting on the
# date-like and digit-like cells *before* assigning their values so Excel's
# automatic type inference doesn't turn them into a real date serial /
# number (which would change their stored type and displayed value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 70

$ws.Range("A$row").NumberFormat = "@"
$ws.Range("C$row").NumberFormat = "@"

$ws.Range("A$row").Value = "2025-11-25"
$ws.Range("B$row").Value = "Pick 4"
$ws.Range("C$row").Value = "251125"
$ws.Range("D$row").Value = "5-4-0-7"
$ws.Range("E$row").Value = "2025-11-25T21:41:29.837+04:00"
